# Add an "Address" column to the teachers list.
#
# A new column is inserted immediately before the existing "District"
# column (old column F), which shifts "District" one column to the right
# (old F -> new G) and gives us a blank column F to fill in as "Address".
#
# The address text for each teacher is pulled out of the second line of
# their "NAMES AND ADDRESS" cell (column B), which reads
#   "<name>
#    <school/address>, <taluk>, <district>."
# The trailing district segment is dropped (it's already captured by the
# "District" column) and the remaining comma-separated segments are
# concatenated together with no separator, matching the target data.
# Rows whose column-B text doesn't follow the two-line "name / address"
# pattern (e.g. continuation rows) are left blank, same as in the source.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("F").Insert()
$ws.Range("F2").Value = "Address"

for ($r = 3; $r -le 42; $r++) {
    $b = $ws.Cells.Item($r, 2).Value2
    if ($b -eq $null) { continue }

    $lines = $b -split "`n"
    if ($lines.Length -ne 2) { continue }

    $parts = $lines[1] -split ","
    $segments = @()
    foreach ($p in $parts) {
        $t = $p.Trim()
        if ($t -ne "") { $segments += $t }
    }

    if ($segments.Length -lt 2) { continue }

    $address = [string]::Join("", $segments[0..($segments.Length - 2)])
    $ws.Cells.Item($r, 6).Value = $address
}
